$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture row 8's existing string values (name/abbreviation/code/parent_name)
# before we touch anything, so we can duplicate them into the new row 9.
$a8 = $ws.Range("A8").Value2
$b8 = $ws.Range("B8").Value2
$c8 = $ws.Range("C8").Value2
$d8 = $ws.Range("D8").Value2

# Row 8 had its warning/danger threshold values swapped:
#   E8 (warning_threshold): 15.0 -> 16.0
#   F8 (danger_threshold):  16.0 -> 15.0
$ws.Range("E8").Value = 16.0
$ws.Range("F8").Value = 15.0

# Add a brand-new row 9 that duplicates row 8's text columns and carries the
# corrected 16.0 / 15.0 threshold values.
$ws.Range("A9").Value = $a8
$ws.Range("B9").Value = $b8
$ws.Range("C9").Value = $c8
$ws.Range("D9").Value = $d8
$ws.Range("E9").Value = 16.0
$ws.Range("F9").Value = 15.0

# Match row 9's cell formatting/style to row 8 so it looks identical.
$ws.Range("A8:F8").Copy()
$ws.Range("A9:F9").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
